$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H129").Value = 1348.9491
$ws.Range("I129").Value = 446.7
$ws.Range("J129").Value = 1811.641
$ws.Range("K129").Value = 1340.1
$ws.Range("L129").Value = 5434.923000000001
$ws.Range("M129").Value = 3659.9
$ws.Range("N129").Value = -15434.923
$ws.Range("H137").Value = 1228.4822
$ws.Range("I137").Value = 1061.7
$ws.Range("J137").Value = 1645.4375
$ws.Range("K137").Value = 3185.1
$ws.Range("L137").Value = 4936.3125
$ws.Range("M137").Value = -635.1000000000004
$ws.Range("N137").Value = -10036.3125
$ws.Range("H138").Value = 1604.47
$ws.Range("I138").Value = 656.9524
$ws.Range("J138").Value = 2290.6035
$ws.Range("K138").Value = 1970.8572
$ws.Range("L138").Value = 6871.810500000001
$ws.Range("M138").Value = 3169.1428
$ws.Range("N138").Value = -17151.8105

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 926432.5600000001
$ws.Range("I32").Value = 1102029.1
$ws.Range("K32").Value = 1102029.1
$ws.Range("M32").Value = -1101742.1
$ws.Range("H61").Value = 2390.8071
$ws.Range("I61").Value = 2140.2163
$ws.Range("J61").Value = 2854.4
$ws.Range("K61").Value = 2140.2163
$ws.Range("L61").Value = 2854.4
$ws.Range("M61").Value = -1928.2163
$ws.Range("N61").Value = -3278.4
$ws.Range("H136").Value = 2390.8071
$ws.Range("I136").Value = 2140.2163
$ws.Range("J136").Value = 2854.4
$ws.Range("K136").Value = 6420.6489
$ws.Range("L136").Value = 8563.200000000001
$ws.Range("M136").Value = -3870.6489
$ws.Range("N136").Value = -13663.2
$ws.Range("H137").Value = 22254.166
$ws.Range("J137").Value = 20210.295
$ws.Range("L137").Value = 20210.295
$ws.Range("N137").Value = -30410.295
$ws.Range("H140").Value = 92696.75
$ws.Range("J140").Value = 92696.75
$ws.Range("L140").Value = 92696.75
$ws.Range("N140").Value = -103056.75

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H105").Value = 62500450
$ws.Range("I105").Value = 62500450
$ws.Range("J105").Value = 0
$ws.Range("K105").Value = 62500450
$ws.Range("L105").Value = 0
$ws.Range("M105").Value = -62498703
$ws.Range("H128").Value = 3496.6667
$ws.Range("I128").Value = 3496.6667
$ws.Range("K128").Value = 10490.0001
$ws.Range("M128").Value = -8000.000100000001
$ws.Range("N105").ClearContents()

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 6791.7856
$ws.Range("I31").Value = 2462.8572
$ws.Range("J31").Value = 8234.762000000001
$ws.Range("K31").Value = 2462.8572
$ws.Range("L31").Value = 8234.762000000001
$ws.Range("M31").Value = -2167.8572
$ws.Range("N31").Value = -8824.762000000001
$ws.Range("H34").Value = 6791.7856
$ws.Range("I34").Value = 2462.8572
$ws.Range("J34").Value = 8234.762000000001
$ws.Range("K34").Value = 2462.8572
$ws.Range("L34").Value = 8234.762000000001
$ws.Range("M34").Value = -2260.8572
$ws.Range("N34").Value = -8638.762000000001
$ws.Range("H58").Value = 1184.225
$ws.Range("I58").Value = 910.5
$ws.Range("J58").Value = 1518.7778
$ws.Range("K58").Value = 910.5
$ws.Range("L58").Value = 1518.7778
$ws.Range("M58").Value = -707.5
$ws.Range("N58").Value = -1924.7778
$ws.Range("H86").Value = 2213.4285
$ws.Range("I86").Value = 2280.2727
$ws.Range("K86").Value = 2280.2727
$ws.Range("M86").Value = -1157.2727
$ws.Range("H89").Value = 2213.4285
$ws.Range("I89").Value = 2280.2727
$ws.Range("K89").Value = 11401.3635
$ws.Range("M89").Value = -5785.363499999999
$ws.Range("H136").Value = 1184.225
$ws.Range("I136").Value = 910.5
$ws.Range("J136").Value = 1518.7778
$ws.Range("K136").Value = 2731.5
$ws.Range("L136").Value = 4556.3334
$ws.Range("M136").Value = -181.5
$ws.Range("N136").Value = -9656.3334

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H2").Value = 46.275864
$ws.Range("I2").Value = 6.3333335
$ws.Range("J2").Value = 89.07143000000001
$ws.Range("K2").Value = 38.000001
$ws.Range("L2").Value = 534.42858
$ws.Range("M2").Value = 74.999999
$ws.Range("N2").Value = -760.42858
$ws.Range("H113").Value = 709.093
$ws.Range("I113").Value = 632.1539
$ws.Range("J113").Value = 826.7646999999999
$ws.Range("K113").Value = 1896.4617
$ws.Range("L113").Value = 2480.2941
$ws.Range("M113").Value = 273.5382999999999
$ws.Range("N113").Value = -6820.2941
$ws.Range("H119").Value = 3399.8
$ws.Range("I119").Value = 2666.3333
$ws.Range("J119").Value = 4500
$ws.Range("K119").Value = 7998.999899999999
$ws.Range("L119").Value = 13500
$ws.Range("M119").Value = -3160.999899999999
$ws.Range("N119").Value = -23176
$ws.Range("H120").Value = 10799.8
$ws.Range("I120").Value = 4999
$ws.Range("J120").Value = 12250
$ws.Range("K120").Value = 14997
$ws.Range("L120").Value = 36750
$ws.Range("M120").Value = -10159
$ws.Range("N120").Value = -46426
$ws.Range("H121").Value = 1230.1
$ws.Range("I121").Value = 470
$ws.Range("J121").Value = 1461.4348
$ws.Range("K121").Value = 1410
$ws.Range("L121").Value = 4384.3044
$ws.Range("M121").Value = -100
$ws.Range("N121").Value = -7004.3044

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H2").Value = 28.5
$ws.Range("I2").Value = 31.8
$ws.Range("J2").Value = 12
$ws.Range("K2").Value = 31.8
$ws.Range("L2").Value = 12
$ws.Range("M2").Value = 81.2
$ws.Range("N2").Value = -238
$ws.Range("H70").Value = 9824.190000000001
$ws.Range("I70").Value = 10963
$ws.Range("J70").Value = 6180
$ws.Range("K70").Value = 10963
$ws.Range("L70").Value = 6180
$ws.Range("M70").Value = -10693
$ws.Range("N70").Value = -6720
$ws.Range("H73").Value = 9824.190000000001
$ws.Range("I73").Value = 10963
$ws.Range("J73").Value = 6180
$ws.Range("K73").Value = 10963
$ws.Range("L73").Value = 6180
$ws.Range("M73").Value = -10027
$ws.Range("N73").Value = -8052
$ws.Range("H132").Value = 3051.5
$ws.Range("I132").Value = 2698
$ws.Range("J132").Value = 3304
$ws.Range("K132").Value = 8094
$ws.Range("L132").Value = 9912
$ws.Range("M132").Value = -5564
$ws.Range("N132").Value = -14972

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H35").Value = 1196
$ws.Range("I35").Value = 1196
$ws.Range("J35").Value = 0
$ws.Range("K35").Value = 1196
$ws.Range("L35").Value = 0
$ws.Range("M35").Value = -860
$ws.Range("H136").Value = 16669197
$ws.Range("I136").Value = 3667
$ws.Range("J136").Value = 27779550
$ws.Range("K136").Value = 11001
$ws.Range("L136").Value = 83338650
$ws.Range("M136").Value = -8451
$ws.Range("N136").Value = -83343750
$ws.Range("N35").ClearContents()

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H123").Value = 29738.166
$ws.Range("J123").Value = 49214.5
$ws.Range("L123").Value = 49214.5
$ws.Range("N123").Value = -59014.5
$ws.Range("H132").Value = 2165192.5
$ws.Range("I132").Value = 577.2586
$ws.Range("J132").Value = 8772965
$ws.Range("K132").Value = 1731.7758
$ws.Range("L132").Value = 26318895
$ws.Range("M132").Value = 798.2242000000001
$ws.Range("N132").Value = -26323955
$ws.Range("H136").Value = 2077.1184
$ws.Range("I136").Value = 1735.9517
$ws.Range("K136").Value = 5207.855100000001
$ws.Range("M136").Value = -2657.855100000001
